$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "culture" column (B) held the value "cholerae" for two samples;
# update those cells to the more precise "V. cholerae" label.
$ws.Range("B10").Value = "V. cholerae"
$ws.Range("B12").Value = "V. cholerae"

# Restore the last active selection on the sheet.
$ws.Range("O20").Select()
